$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage, without leaving a
# lasting style/number-format change on the cell (price strings like
# "211.45" would otherwise be auto-coerced to a float by Excel).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.710.63"
$ws.Range("D3").Value = "1.600.31"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.24%  "
Set-TextValue "D5" "211.45"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +0.90%  "
Set-TextValue "D11" "0.0843"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.824.72"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "1.601.91"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +0.49%  "
Set-TextValue "D16" "65.36"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "26.686.40"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("E18").Value = "  +3.65%  "
$ws.Range("E19").Value = "  +0.15%  "
Set-TextValue "D20" "209.47"
$ws.Range("E20").Value = "  +0.58%  "
Set-TextValue "D21" "7.19"
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +0.92%  "
Set-TextValue "D25" "142.87"
$ws.Range("E25").Value = "  -1.56%  "
Set-TextValue "D26" "1.01"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("E28").Value = "  +0.20%  "
Set-TextValue "D29" "15.35"
$ws.Range("E29").Value = "  +0.68%  "
Set-TextValue "D30" "0.0518"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "1.292.60"
$ws.Range("E34").Value = "  +1.06%  "
Set-TextValue "D35" "0.620"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +0.13%  "
Set-TextValue "D39" "1.08"
$ws.Range("E39").Value = "  +17.96%  "
Set-TextValue "D40" "0.826"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -0.38%  "
Set-TextValue "D44" "63.15"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").Value = "1.735.97"
$ws.Range("E45").Value = "  +0.45%  "
Set-TextValue "D46" "90.92"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("E48").Value = "  -1.08%  "
Set-TextValue "D49" "0.0510"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("E50").Value = "  +0.15%  "
Set-TextValue "D51" "7.36"
$ws.Range("E51").Value = "  -0.87%  "
